$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the NroSiniestro (column E) claim numbers, preserving the
# "number stored as text" formatting (leading apostrophe keeps the
# existing quotePrefix style and exact whitespace) for each row.
$ws.Range("E2").Value = "'1120194100448  "
$ws.Range("E3").Value = "'1220194200694 "
$ws.Range("E5").Value = "'0420172010228    "
$ws.Range("E6").Value = "'1220170301466    "
$ws.Range("E7").Value = "'1120170200969    "

# Swap the tester for the fifth data row (preproducciongestion env).
$ws.Range("C5").Value = "eceldane"

# Add a new row (row 8) with a stray Usuario value, matching the
# corrected test case behavior.
$ws.Range("C8").Value = "mpimpignano"

# Update the active cell selection left after the edits.
$null = $ws.Range("F4").Select()
